$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.825.09"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.499.84"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'598.47"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'194.33"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.208"
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").Value = "'0.651"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "'53.67"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'0.0000300"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").Value = "'9.51"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "4.049.33"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "'608.08"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "69.916.91"
$ws.Range("D17").Value = "'18.97"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'12.57"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "3.489.88"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "'0.990"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'18.02"
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").Value = "'104.52"
$ws.Range("E23").Value = "  +9.16%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").Value = "'4.56"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").Value = "'3.06"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'9.71"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "'33.57"
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").Value = "'4.60"
$ws.Range("E30").Value = "  +26.70%  "
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'12.62"
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").Value = "'0.115"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "3.735.44"
$ws.Range("E35").Value = "  +5.93%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0808"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'3.05"
$ws.Range("E38").Value = "  -6.03%  "
$ws.Range("D39").Value = "'0.391"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.57"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'36.56"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'502.30"
$ws.Range("E42").Value = "  -6.45%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'0.0457"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "'3.32"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'2.82"
$ws.Range("E47").Value = "  -3.99%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'8.73"
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("D50").Value = "'131.75"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'0.000242"
$ws.Range("E51").Value = "  -0.31%  "
